$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 20.91598510742188
$ws.Range("C3").Value = 18.91684532165527
$ws.Range("C4").Value = 18.49102973937988
$ws.Range("C5").Value = 19.8667049407959
$ws.Range("C6").Value = 21.57688140869141
